$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 16:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1647212
$ws.Range("C4").Value = 2118
$ws.Range("E4").Value = 1146204
$ws.Range("G4").Value = 49
$ws.Range("H4").Value = 97696

# Row 6 - Brasil
$ws.Range("B6").Value = 334777
$ws.Range("C6").Value = 3887
$ws.Range("E6").Value = 178132
$ws.Range("G6").Value = 167
$ws.Range("H6").Value = 21215

# Row 60 - Oman
$ws.Range("E60").Value = 5374
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 35

# Row 137 - Reunion
$ws.Range("B137").Value = 452
$ws.Range("C137").Value = 3
$ws.Range("E137").Value = 40

# Row 139 - Estado de Palestina
$ws.Range("E139").Value = 74
$ws.Range("G139").Value = 1
$ws.Range("H139").Value = 3
